# Rever_DailyTrack_BALRAJ_2022.xlsx - "Add files via upload"
# Fill in the two new daily-task rows (row 25 and row 26) on the FEB-22
# sheet, and move the active-cell selection to F26 (the last cell typed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")

# --- Row 25: a brand-new entry (No. 18, dated 22-Feb-2022) -----------------
# Copy the date format from the row above (B24) onto B25 so the new date
# cell picks up the existing "m/d/yyyy" style instead of creating a new one.
$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)

$ws.Range("A25").Value = 18
$ws.Range("B25").Value = 44614
$ws.Range("C25").Value = "RPA GSS"
$ws.Range("D25").Value = "1. The updating master data file has been completed 50% and other process is work in progress"
$ws.Range("E25").Value = 0.7
$ws.Range("F25").Value = "WIP"

# --- Row 26: continuation line for the same entry (only D/E/F filled) -----
$ws.Range("D26").Value = "2. Getting cell value along with column number is work in progress"
$ws.Range("E26").Value = 0.7
$ws.Range("F26").Value = "WIP"

# --- Move the selection to the last edited cell, matching the saved view --
$ws.Range("F26").Select()
